# PM15 Tidsregistrering for Laila.xlsx - add three new time-tracking entries
# (rows 42-44 on "Ark1") and extend the "Timer i alt" shared formula to
# cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 42: Systemtest OC0802 tjekke eksempler / Tester -------------------
$ws.Range("A42").Value = "Systemtest OC0802 tjekke eksempler"
$ws.Range("B42").Value = "Tester"
$ws.Range("C42").Value = 43901
$ws.Range("D42").Value = 0.3611111111111111
$ws.Range("E42").Value = 0.39930555555555558

# --- Row 43: Review af AD09 med Tommy / Any Role ----------------------------
$ws.Range("A43").Value = "Review af AD09 med Tommy"
$ws.Range("B43").Value = "Any Role"
$ws.Range("C43").Value = 43901
$ws.Range("D43").Value = 0.39930555555555558
$ws.Range("E43").Value = 0.40625

# --- Row 44: Systemtest OC0802 / Tester -------------------------------------
$ws.Range("A44").Value = "Systemtest OC0802"
$ws.Range("B44").Value = "Tester"
$ws.Range("C44").Value = 43901
$ws.Range("D44").Value = 0.41319444444444442
$ws.Range("E44").Value = 0.63194444444444442

# Fill in "Timer i alt" (G) for the new rows, matching the existing
# E-D shared formula used throughout the column.
$ws.Range("G42:G44").Formula = "=E42-D42"

# Update the active selection to match the latest edit location.
$ws.Range("F44").Select() | Out-Null
